$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.87621524227613157
$ws.Range("BO1").Value = 0.7436005708065434
$ws.Range("D2").Value = 0.91854721851115095
$ws.Range("AW2").Value = 0.50962223072270207
$ws.Range("A3").Value = 0.69081635267734098
$ws.Range("B3").Value = 0.72370018584838192
$ws.Range("P3").Value = 0.91797083768054888
$ws.Range("BN3").Value = 0.84271467483268081
$ws.Range("C4").Value = 0.84217198647965674
$ws.Range("C5").Value = 0.99188293513677406
$ws.Range("D5").Value = 0.88260552884630761
$ws.Range("G5").Value = 0.99823415198058285
$ws.Range("D6").Value = 0.82299740129551302
$ws.Range("H6").Value = 0.67948575165955094
$ws.Range("F7").Value = 0.89066897884417662
$ws.Range("W7").Value = 0.88875698598024355
$ws.Range("AP7").Value = 0.90244496549211428
$ws.Range("J8").Value = 0.60523639465998369
$ws.Range("G9").Value = 0.76751339466123825
$ws.Range("H9").Value = 0.76191949308465878
$ws.Range("S9").Value = 0.83768495006770771
$ws.Range("J11").Value = 0.95196568966429918
$ws.Range("L11").Value = 0.71621165273142751
$ws.Range("M11").Value = 0.76508565488746161
$ws.Range("G12").Value = 0.65527927952799447
$ws.Range("J12").Value = 0.65846519606561327
$ws.Range("M12").Value = 0.98930680414585515
$ws.Range("N13").Value = 0.57196236061590044
$ws.Range("AD13").Value = 0.78030777396791007
$ws.Range("M15").Value = 0.75900808059269342
$ws.Range("N15").Value = 0.91841520740187099
$ws.Range("P15").Value = 0.74664986131569866
$ws.Range("I16").Value = 0.85156658165090304
$ws.Range("AP16").Value = 0.78942114948224196
$ws.Range("O17").Value = 0.93250241593383987
$ws.Range("R17").Value = 0.91094512305024367
$ws.Range("BM17").Value = 0.86573759213312274
$ws.Range("AT18").Value = 0.91703726775060113
$ws.Range("T19").Value = 0.89976660192142277
$ws.Range("AW19").Value = 0.70576554110281675
$ws.Range("BO20").Value = 0.93069251392328223
$ws.Range("S21").Value = 0.78462576227877201
$ws.Range("V21").Value = 0.8403627235731771
$ws.Range("T22").Value = 0.86040872515730848
$ws.Range("X22").Value = 0.9154750315600898
$ws.Range("AF22").Value = 0.99494818900374404
$ws.Range("U23").Value = 0.77731400787962457
$ws.Range("Y23").Value = 0.88675126490195932
$ws.Range("Z25").Value = 0.79189637583180128
$ws.Range("AA25").Value = 0.97700152208141244
$ws.Range("X26").Value = 0.92673949229158625
$ws.Range("AA26").Value = 0.94505646047355518
$ws.Range("AB26").Value = 0.8648102811807814
$ws.Range("AB27").Value = 0.77349969488731674
$ws.Range("AC27").Value = 0.95391436547275443
$ws.Range("AD28").Value = 0.67402887820914636
$ws.Range("AB29").Value = 0.87493081860372068
$ws.Range("AF30").Value = 0.92228449794005818
$ws.Range("W31").Value = 0.56943724405029394
$ws.Range("AC31").Value = 0.93947804491013454
$ws.Range("AG31").Value = 0.69178177792789675
$ws.Range("AG32").Value = 0.78865999308167511
$ws.Range("C33").Value = 0.6753458585654637
$ws.Range("X34").Value = 0.83281880618645232
$ws.Range("AF34").Value = 0.79078091550332519
$ws.Range("AJ34").Value = 0.94597970713851653
$ws.Range("AG35").Value = 0.84580000915902398
$ws.Range("AH35").Value = 0.91563600586571836
$ws.Range("AK35").Value = 0.9692529809201782
$ws.Range("BB35").Value = 0.75453591669058051
$ws.Range("Z36").Value = 0.92676477038058647
$ws.Range("AL37").Value = 0.86849976901782378
$ws.Range("AE38").Value = 0.72766061045847497
$ws.Range("AJ38").Value = 0.79932680874835893
$ws.Range("AN38").Value = 0.99527770560479223
$ws.Range("AK39").Value = 0.91357773193732494
$ws.Range("BA39").Value = 0.7378483169503649
$ws.Range("AA40").Value = 0.89293915307066163
$ws.Range("AM40").Value = 0.98076145324385944
$ws.Range("AO40").Value = 0.76518962079030894
$ws.Range("AM41").Value = 0.80759772547781128
$ws.Range("AP41").Value = 0.99989019617003683
$ws.Range("AQ41").Value = 0.72606658442228633
$ws.Range("AQ42").Value = 0.90418818135281454
$ws.Range("BH43").Value = 0.80762626388257674
$ws.Range("BP43").Value = 0.97368038295869419
$ws.Range("AP44").Value = 0.76937168708756043
$ws.Range("AQ44").Value = 0.99050802042801867
$ws.Range("AQ45").Value = 0.97569954179011575
$ws.Range("AR45").Value = 0.60724512748396442
$ws.Range("AT45").Value = 0.79023144502540421
$ws.Range("AU45").Value = 0.99326925737043337
$ws.Range("E46").Value = 0.97624751982568048
$ws.Range("AT47").Value = 0.91647924999173302
$ws.Range("AV47").Value = 0.96307062895154527
$ws.Range("BB47").Value = 0.91663335544095559
$ws.Range("AT48").Value = 0.69940838303915542
$ws.Range("AW48").Value = 0.96576586239824347
$ws.Range("AZ49").Value = 0.75354920782961132
$ws.Range("AV50").Value = 0.72200084993616054
$ws.Range("R51").Value = 0.7818640068799223
$ws.Range("AX51").Value = 0.87879337432650662
$ws.Range("AZ51").Value = 0.79127609549480404
$ws.Range("AX52").Value = 0.94053067889379827
$ws.Range("N53").Value = 0.8564133346268602
$ws.Range("AY53").Value = 0.73173508712314805
$ws.Range("BC53").Value = 0.71303807507401373
$ws.Range("AC54").Value = 0.91467713430775432
$ws.Range("AZ54").Value = 0.8882501732559569
$ws.Range("BC54").Value = 0.72021331589158033
$ws.Range("BD54").Value = 0.86119446234941144
$ws.Range("BC56").Value = 0.57720223409323301
$ws.Range("BC57").Value = 0.72797300709090229
$ws.Range("BD57").Value = 0.58348215073425935
$ws.Range("BD58").Value = 0.87958976100973629
$ws.Range("BE59").Value = 0.89250401560291981
$ws.Range("BF59").Value = 0.62108269205191668
$ws.Range("BF60").Value = 0.98122392554319071
$ws.Range("BG60").Value = 0.78424250675863261
$ws.Range("BG61").Value = 0.94412030023282256
$ws.Range("BK61").Value = 0.6723062605687391
$ws.Range("AY62").Value = 0.8491148977069165
$ws.Range("BH62").Value = 0.65447015082025184
$ws.Range("BI62").Value = 0.96835771742006505
$ws.Range("BL62").Value = 0.65841413690698136
$ws.Range("AQ63").Value = 0.76862998326669141
$ws.Range("BM63").Value = 0.9595568334129142
$ws.Range("Q64").Value = 0.99529305351599118
$ws.Range("BM64").Value = 0.93028987234634242
$ws.Range("BN65").Value = 0.92840062839388149
$ws.Range("AR66").Value = 0.80163830353965704
$ws.Range("BL66").Value = 0.61504158447522395
$ws.Range("BO66").Value = 0.64990144009219319
$ws.Range("BE67").Value = 0.74315842911938135
$ws.Range("AN68").Value = 0.994827999712671
$ws.Range("BO68").Value = 0.62514581553151749
